# "10 years Finalization data"
# Duplicate the daily-data table (header row 9 + data rows 10-40 of the
# original "Data Harian - Table" sheet) into a brand-new "Sheet1" worksheet,
# re-based so the table starts at A1. The new sheet becomes the active /
# selected tab, matching a Finalized/clean extract of the same month's data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet, placed right after the source sheet; Excel names it "Sheet1".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Copy the table (header + 31 daily rows) including formatting, then values.
$src = $ws1.Range("A9:K40")
$src.Copy()

$dst = $ws2.Range("A1")
$dst.PasteSpecial(-4122)   # xlPasteFormats
$dst.PasteSpecial(-4163)   # xlPasteValues

# Match the on-screen selection state from the edit.
$ws2.Range("A1:K32").Select()

$ws1.Activate()
$ws1.Range("A9:K40").Select()

$ws2.Activate()
